$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New template-type rows for the OTP email content/subject templates, one
# pair per existing language (eng, fra, ara) already present in the table
# (columns: lang_code | code | descr | is_active). Insert each new row as a
# duplicate of the last existing data row (442) so it inherits the same
# cell styles (language-code style + boolean/text style on is_active),
# then overwrite the cell values.
for ($i = 0; $i -lt 6; $i++) {
    $ws.Rows(442).Copy()
    $ws.Rows(443).Insert(-4121)
}

$newRows = @(
    @("eng", "otp-email-content-template", "Template for OTP Email Content"),
    @("eng", "otp-email-subject-template", "Template for OTP Email Subject"),
    @("fra", "otp-email-content-template", "Template for OTP Email Content"),
    @("fra", "otp-email-subject-template", "Template for OTP Email Subject"),
    @("ara", "otp-email-content-template", "Template for OTP Email Content"),
    @("ara", "otp-email-subject-template", "Template for OTP Email Subject")
)

# Write column A (lang codes - all already existing shared strings) for
# every new row first.
$startRow = 443
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value = $newRows[$i][0]
}

# Write column B ("code") only for the first eng pair so the two brand new
# strings ("otp-email-content-template" / "otp-email-subject-template")
# are appended to the shared-string table before any column C strings.
$ws.Range("B443").Value = $newRows[0][1]
$ws.Range("B444").Value = $newRows[1][1]

# Write column C ("descr") for the first eng pair next, appending the two
# new description strings right after the two new code strings.
$ws.Range("C443").Value = $newRows[0][2]
$ws.Range("C444").Value = $newRows[1][2]

# Remaining rows (fra/ara pairs) only reuse already-appended shared
# strings, so order no longer matters.
for ($i = 2; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
}

# Keep the view pointed near the bottom of the table, mirroring the
# author's scroll position after adding the rows.
$ws.Application.Goto($ws.Range("A428"))
$ws.Range("C445").Select()
